# "Actualizacion de tarea completada"
# - Task in row 27 (C27) flips from the "en proceso" text flag to a
#   checked/complete flag (numeric 1, percent-formatted), matching the
#   other completed rows (e.g. row 26).
# - A new task row is appended (row 29) plus a trailing formatted blank
#   row (row 30), and the view selection moves to B4 with the sticky
#   top-left-cell scroll position cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 27: mark task complete - replace the "en proceso" label with a
# percent-formatted 1 (same formatting family as C26/C3/etc.).
$ws.Range("C27").NumberFormat = "0%"
$ws.Range("C27").Value = 1

# New task row.
$ws.Range("A29").Value = "Validacion de cuit para mostrar mensaje correcto"

# Trailing blank separator row, underlined like the other section breaks
# (C18/C13/D11).
$ws.Range("A30").Font.Underline = 2

# Update the view: clear the scrolled top-left cell and move the
# selection to B4.
$ws.Activate()
$ws.Range("B4").Select() | Out-Null
